# Updates implied-volatility re-calculation results (lastPrice, bid, ask, change,
# percentChange, volume, impliedVolatility, Delta, Gamma, Vega, Rho, Theta, inTheMoney)
# for rows 4-31 of the puts options chain sheet, per the refreshed pricing run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("L4").Value = 0.42188078125
# Row 5
$ws.Range("L5").Value = 0.335944140625
# Row 6
$ws.Range("E6").Value = 0.02
$ws.Range("F6").Value = 0.02
$ws.Range("G6").Value = 0.03
$ws.Range("H6").Value = -0.02
$ws.Range("I6").Value = -50
$ws.Range("J6").Value = 1067
$ws.Range("L6").Value = 0.22852333984375
$ws.Range("P6").Value = -0.001
$ws.Range("R6").Value = 0.001
$ws.Range("T6").Value = -0.006
# Row 7
$ws.Range("F7").Value = 0.02
$ws.Range("G7").Value = 0.03
$ws.Range("J7").Value = 1776
$ws.Range("L7").Value = 0.179695703125
$ws.Range("R7").Value = 0.001
$ws.Range("T7").Value = -0.007
# Row 8
$ws.Range("F8").Value = 0.03
$ws.Range("G8").Value = 0.04
$ws.Range("J8").Value = 3637
$ws.Range("L8").Value = 0.1601646484375
$ws.Range("P8").Value = -0.002
$ws.Range("Q8").Value = 0.002
$ws.Range("T8").Value = -0.011
# Row 9
$ws.Range("J9").Value = 7571
$ws.Range("L9").Value = 0.133797724609375
$ws.Range("T9").Value = -0.012
# Row 10
$ws.Range("E10").Value = 0.04
$ws.Range("H10").Value = -0.07000000000000001
$ws.Range("I10").Value = -63.636364
$ws.Range("J10").Value = 3823
$ws.Range("L10").Value = 0.12500875
$ws.Range("P10").Value = -0.004
$ws.Range("Q10").Value = 0.004
# Row 11
$ws.Range("E11").Value = 0.05
$ws.Range("F11").Value = 0.04
$ws.Range("G11").Value = 0.05
$ws.Range("H11").Value = -0.09
$ws.Range("I11").Value = -64.28572
$ws.Range("J11").Value = 25580
$ws.Range("L11").Value = 0.110360458984375
$ws.Range("P11").Value = -0.005
$ws.Range("Q11").Value = 0.006
$ws.Range("R11").Value = 0.003
$ws.Range("T11").Value = -0.018
# Row 12
$ws.Range("E12").Value = 0.06
$ws.Range("F12").Value = 0.06
$ws.Range("G12").Value = 0.07000000000000001
$ws.Range("H12").Value = -0.13
$ws.Range("I12").Value = -68.42104999999999
$ws.Range("J12").Value = 28590
$ws.Range("L12").Value = 0.102548037109375
$ws.Range("P12").Value = -0.01
$ws.Range("Q12").Value = 0.011
$ws.Range("R12").Value = 0.006
$ws.Range("T12").Value = -0.029
# Row 13
$ws.Range("E13").Value = 0.09
$ws.Range("F13").Value = 0.08
$ws.Range("G13").Value = 0.09
$ws.Range("H13").Value = -0.19
$ws.Range("I13").Value = -67.85714
$ws.Range("J13").Value = 42621
$ws.Range("L13").Value = 0.09229423339843749
$ws.Range("P13").Value = -0.017
$ws.Range("Q13").Value = 0.02
$ws.Range("R13").Value = 0.01
$ws.Range("T13").Value = -0.043
# Row 14
$ws.Range("E14").Value = 0.15
$ws.Range("F14").Value = 0.14
$ws.Range("G14").Value = 0.15
$ws.Range("H14").Value = -0.29
$ws.Range("I14").Value = -65.90909000000001
$ws.Range("J14").Value = 62342
$ws.Range("L14").Value = 0.08643491699218751
$ws.Range("P14").Value = -0.039
$ws.Range("Q14").Value = 0.043
$ws.Range("R14").Value = 0.019
$ws.Range("S14").Value = -0
$ws.Range("T14").Value = -0.082
# Row 15
$ws.Range("E15").Value = 0.25
$ws.Range("F15").Value = 0.25
$ws.Range("G15").Value = 0.26
$ws.Range("H15").Value = -0.41000003
$ws.Range("I15").Value = -62.121212
$ws.Range("J15").Value = 101328
$ws.Range("L15").Value = 0.08155215332031251
$ws.Range("P15").Value = -0.093
$ws.Range("Q15").Value = 0.089
$ws.Range("R15").Value = 0.038
$ws.Range("S15").Value = -0.001
$ws.Range("T15").Value = -0.149
# Row 16
$ws.Range("E16").Value = 0.46
$ws.Range("F16").Value = 0.44
$ws.Range("G16").Value = 0.45
$ws.Range("H16").Value = -0.52
$ws.Range("I16").Value = -53.061222
$ws.Range("J16").Value = 188487
$ws.Range("L16").Value = 0.07727973510742189
$ws.Range("P16").Value = -0.204
$ws.Range("Q16").Value = 0.161
$ws.Range("R16").Value = 0.065
$ws.Range("S16").Value = -0.002
$ws.Range("T16").Value = -0.237
# Row 17
$ws.Range("E17").Value = 0.75
$ws.Range("F17").Value = 0.74
$ws.Range("G17").Value = 0.75
$ws.Range("H17").Value = -0.65999997
$ws.Range("I17").Value = -46.80851
$ws.Range("J17").Value = 99968
$ws.Range("L17").Value = 0.07239697143554688
$ws.Range("M17").Value = $False
$ws.Range("P17").Value = -0.391
$ws.Range("Q17").Value = 0.232
$ws.Range("R17").Value = 0.08799999999999999
$ws.Range("S17").Value = -0.005
$ws.Range("T17").Value = -0.293
# Row 18
$ws.Range("E18").Value = 1.18
$ws.Range("F18").Value = 1.17
$ws.Range("G18").Value = 1.18
$ws.Range("H18").Value = -0.7900001
$ws.Range("I18").Value = -40.101524
$ws.Range("J18").Value = 38911
$ws.Range("L18").Value = 0.06519489501953124
$ws.Range("P18").Value = -0.643
$ws.Range("Q18").Value = 0.251
$ws.Range("R18").Value = 0.08500000000000001
$ws.Range("T18").Value = -0.237
# Row 19
$ws.Range("E19").Value = 1.75
$ws.Range("F19").Value = 1.72
$ws.Range("G19").Value = 1.75
$ws.Range("H19").Value = -0.9000001
$ws.Range("I19").Value = -33.96227
$ws.Range("J19").Value = 7373
$ws.Range("L19").Value = 0.05054660400390625
$ws.Range("P19").Value = -0.909
$ws.Range("Q19").Value = 0.142
$ws.Range("R19").Value = 0.037
$ws.Range("S19").Value = -0.011
$ws.Range("T19").Value = -0.037
# Row 20
$ws.Range("E20").Value = 2.5
$ws.Range("F20").Value = 2.4
$ws.Range("G20").Value = 2.5
$ws.Range("H20").Value = -0.9100001
$ws.Range("I20").Value = -26.686218
$ws.Range("J20").Value = 1882
$ws.Range("L20").Value = 0.00001
$ws.Range("P20").Value = -1
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 0
$ws.Range("S20").Value = -0.012
$ws.Range("T20").Value = 0.063
# Row 21
$ws.Range("E21").Value = 3.34
$ws.Range("F21").Value = 3.03
$ws.Range("G21").Value = 3.55
$ws.Range("H21").Value = -0.9600003
$ws.Range("I21").Value = -22.325588
$ws.Range("J21").Value = 1930
$ws.Range("L21").Value = 0.00001
$ws.Range("P21").Value = -1
$ws.Range("Q21").Value = 0
$ws.Range("R21").Value = 0
$ws.Range("S21").Value = -0.012
$ws.Range("T21").Value = 0.063
# Row 22
$ws.Range("E22").Value = 4.34
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 4.59
$ws.Range("H22").Value = -0.9099998500000001
$ws.Range("I22").Value = -17.33333
$ws.Range("J22").Value = 354
$ws.Range("L22").Value = 0.00001
$ws.Range("P22").Value = -1
$ws.Range("Q22").Value = 0
$ws.Range("R22").Value = 0
$ws.Range("S22").Value = -0.012
$ws.Range("T22").Value = 0.063
# Row 23
$ws.Range("E23").Value = 6.42
$ws.Range("F23").Value = 4.86
$ws.Range("G23").Value = 5.6
$ws.Range("H23").Value = 0.32999992
$ws.Range("I23").Value = 5.418718
$ws.Range("J23").Value = 148
$ws.Range("L23").Value = 0.00001
$ws.Range("P23").Value = -1
$ws.Range("Q23").Value = 0
$ws.Range("R23").Value = 0
$ws.Range("S23").Value = -0.012
$ws.Range("T23").Value = 0.064
# Row 24
$ws.Range("F24").Value = 6
$ws.Range("G24").Value = 6.59
$ws.Range("L24").Value = 0.00001
$ws.Range("P24").Value = -1
$ws.Range("Q24").Value = 0
$ws.Range("R24").Value = 0
$ws.Range("S24").Value = -0.012
$ws.Range("T24").Value = 0.064
# Row 25
$ws.Range("F25").Value = 7
$ws.Range("G25").Value = 7.59
$ws.Range("L25").Value = 0.00001
$ws.Range("P25").Value = -1
$ws.Range("Q25").Value = 0
$ws.Range("R25").Value = 0
$ws.Range("S25").Value = -0.012
$ws.Range("T25").Value = 0.064
# Row 26
$ws.Range("F26").Value = 8
$ws.Range("G26").Value = 8.59
$ws.Range("L26").Value = 0.00001
$ws.Range("P26").Value = -1
$ws.Range("Q26").Value = 0
$ws.Range("R26").Value = 0
$ws.Range("S26").Value = -0.012
$ws.Range("T26").Value = 0.064
# Row 27
$ws.Range("F27").Value = 8.890000000000001
$ws.Range("G27").Value = 9.77
$ws.Range("L27").Value = 0.18946123046875
$ws.Range("P27").Value = -0.985
$ws.Range("Q27").Value = 0.008999999999999999
$ws.Range("R27").Value = 0.008999999999999999
$ws.Range("T27").Value = -0.017
# Row 28
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 10.59
$ws.Range("L28").Value = 0.00001
$ws.Range("P28").Value = -1
$ws.Range("Q28").Value = 0
$ws.Range("R28").Value = 0
$ws.Range("T28").Value = 0.064
# Row 29
$ws.Range("F29").Value = 11
$ws.Range("G29").Value = 11.59
$ws.Range("L29").Value = 0.00001
$ws.Range("P29").Value = -1
$ws.Range("Q29").Value = 0
$ws.Range("R29").Value = 0
$ws.Range("T29").Value = 0.064
# Row 30
$ws.Range("F30").Value = 12.94
$ws.Range("G30").Value = 13.75
$ws.Range("I30").Value = 10.256414
$ws.Range("L30").Value = 0.242195078125
$ws.Range("P30").Value = -0.992
$ws.Range("Q30").Value = 0.004
$ws.Range("R30").Value = 0.005
$ws.Range("T30").Value = 0.003
# Row 31
$ws.Range("F31").Value = 32.78
$ws.Range("G31").Value = 33.81
$ws.Range("L31").Value = 0.5293015820312501
$ws.Range("P31").Value = -0.996
$ws.Range("Q31").Value = 0.001
$ws.Range("R31").Value = 0.003
$ws.Range("T31").Value = -0.003
